$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 92472.875
$ws.Range("I11").Value = 92472.875
$ws.Range("K11").Value = 92472.875
$ws.Range("M11").Value = -92332.875
$ws.Range("H51").Value = 2986.48
$ws.Range("I51").Value = 2985.9167
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 2985.9167
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -2501.9167
$ws.Range("N51").Value = -3968
$ws.Range("H86").Value = 12355.667
$ws.Range("I86").Value = 12360.4
$ws.Range("K86").Value = 12360.4
$ws.Range("M86").Value = -11237.4
$ws.Range("H89").Value = 12355.667
$ws.Range("I89").Value = 12360.4
$ws.Range("K89").Value = 61802
$ws.Range("M89").Value = -56186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 729.38464
$ws.Range("J2").Value = 699.6
$ws.Range("L2").Value = 699.6
$ws.Range("N2").Value = -925.6
$ws.Range("H45").Value = 2575.5454
$ws.Range("J45").Value = 3320
$ws.Range("L45").Value = 3320
$ws.Range("N45").Value = -4074
$ws.Range("H97").Value = 1645.8
$ws.Range("I97").Value = 1257.25
$ws.Range("K97").Value = 1257.25
$ws.Range("M97").Value = -761.25
$ws.Range("H110").Value = 1941.1177
$ws.Range("I110").Value = 1857
$ws.Range("K110").Value = 1857
$ws.Range("M110").Value = 188
$ws.Range("H116").Value = 729.38464
$ws.Range("J116").Value = 699.6
$ws.Range("L116").Value = 699.6
$ws.Range("N116").Value = -5287.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 729.38464
$ws.Range("J3").Value = 699.6
$ws.Range("L3").Value = 699.6
$ws.Range("N3").Value = -927.6
$ws.Range("H134").Value = 7046.1577
$ws.Range("I134").Value = 7430.25
$ws.Range("K134").Value = 22290.75
$ws.Range("M134").Value = -19755.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 148286.58
$ws.Range("I16").Value = 8751.5
$ws.Range("J16").Value = 334333.34
$ws.Range("K16").Value = 8751.5
$ws.Range("L16").Value = 334333.34
$ws.Range("M16").Value = -8464.5
$ws.Range("N16").Value = -334907.34
$ws.Range("H31").Value = 2530.5
$ws.Range("I31").Value = 2530.5
$ws.Range("K31").Value = 2530.5
$ws.Range("M31").Value = -2235.5
$ws.Range("H34").Value = 2530.5
$ws.Range("I34").Value = 2530.5
$ws.Range("K34").Value = 2530.5
$ws.Range("M34").Value = -2328.5
$ws.Range("H41").Value = 10115.308
$ws.Range("I41").Value = 1250.3334
$ws.Range("J41").Value = 17713.857
$ws.Range("K41").Value = 1250.3334
$ws.Range("L41").Value = 17713.857
$ws.Range("M41").Value = -822.3334
$ws.Range("N41").Value = -18569.857
$ws.Range("H62").Value = 4128.2856
$ws.Range("I62").Value = 3779.8
$ws.Range("K62").Value = 3779.8
$ws.Range("M62").Value = -3155.8
$ws.Range("H65").Value = 4128.2856
$ws.Range("I65").Value = 3779.8
$ws.Range("K65").Value = 18899
$ws.Range("M65").Value = -15779
$ws.Range("H86").Value = 91848.73
$ws.Range("I86").Value = 152981.33
$ws.Range("J86").Value = 18489.6
$ws.Range("K86").Value = 152981.33
$ws.Range("L86").Value = 18489.6
$ws.Range("M86").Value = -151858.33
$ws.Range("N86").Value = -20735.6
$ws.Range("H89").Value = 91848.73
$ws.Range("I89").Value = 152981.33
$ws.Range("J89").Value = 18489.6
$ws.Range("K89").Value = 764906.6499999999
$ws.Range("L89").Value = 92448
$ws.Range("M89").Value = -759290.6499999999
$ws.Range("N89").Value = -103680
$ws.Range("H105").Value = 22083.334
$ws.Range("I105").Value = 22083.334
$ws.Range("K105").Value = 22083.334
$ws.Range("M105").Value = -20336.334
$ws.Range("H113").Value = 148286.58
$ws.Range("I113").Value = 8751.5
$ws.Range("J113").Value = 334333.34
$ws.Range("K113").Value = 8751.5
$ws.Range("L113").Value = 334333.34
$ws.Range("M113").Value = -6581.5
$ws.Range("N113").Value = -338673.34
$ws.Range("H122").Value = 9278.370999999999
$ws.Range("I122").Value = 2093.4644
$ws.Range("K122").Value = 6280.3932
$ws.Range("M122").Value = -3830.3932
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1351191.9
$ws.Range("J9").Value = 1501888.8
$ws.Range("L9").Value = 4505666.4
$ws.Range("N9").Value = -4506114.4
$ws.Range("H37").Value = 42499.445
$ws.Range("J37").Value = 42499.445
$ws.Range("L37").Value = 127498.335
$ws.Range("N37").Value = -127722.335
$ws.Range("H101").Value = 8345654
$ws.Range("J101").Value = 8345654
$ws.Range("L101").Value = 25036962
$ws.Range("N101").Value = -25041830
$ws.Range("H121").Value = 15772.056
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 15772.056
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 47316.16800000001
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -49936.16800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3370.1
$ws.Range("I80").Value = 2866.8333
$ws.Range("K80").Value = 2866.8333
$ws.Range("M80").Value = -1868.8333
$ws.Range("H83").Value = 3370.1
$ws.Range("I83").Value = 2866.8333
$ws.Range("K83").Value = 14334.1665
$ws.Range("M83").Value = -9342.166499999999
$ws.Range("H102").Value = 3099.3572
$ws.Range("I102").Value = 3074.25
$ws.Range("K102").Value = 3074.25
$ws.Range("M102").Value = -1452.25
$ws.Range("H113").Value = 1840
$ws.Range("I113").Value = 1857.0435
$ws.Range("K113").Value = 1857.0435
$ws.Range("M113").Value = 312.9565
$ws.Range("H126").Value = 2983.4285
$ws.Range("I126").Value = 2723.5
$ws.Range("K126").Value = 8170.5
$ws.Range("M126").Value = -5700.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1404.4642
$ws.Range("J16").Value = 933.6667
$ws.Range("L16").Value = 933.6667
$ws.Range("N16").Value = -1273.6667
$ws.Range("H55").Value = 1338.2858
$ws.Range("I55").Value = 1142.7333
$ws.Range("K55").Value = 1142.7333
$ws.Range("M55").Value = -969.7333000000001
$ws.Range("H61").Value = 9126.5
$ws.Range("I61").Value = 7639.15
$ws.Range("J61").Value = 24000
$ws.Range("K61").Value = 7639.15
$ws.Range("L61").Value = 24000
$ws.Range("M61").Value = -7437.15
$ws.Range("N61").Value = -24404
$ws.Range("H68").Value = 8606.333000000001
$ws.Range("I68").Value = 8666.416999999999
$ws.Range("J68").Value = 8546.25
$ws.Range("K68").Value = 8666.416999999999
$ws.Range("L68").Value = 8546.25
$ws.Range("M68").Value = -7917.416999999999
$ws.Range("N68").Value = -10044.25
$ws.Range("H71").Value = 8606.333000000001
$ws.Range("I71").Value = 8666.416999999999
$ws.Range("J71").Value = 8546.25
$ws.Range("K71").Value = 43332.085
$ws.Range("L71").Value = 42731.25
$ws.Range("M71").Value = -39588.085
$ws.Range("N71").Value = -50219.25
$ws.Range("H82").Value = 2121
$ws.Range("I82").Value = 2104.4
$ws.Range("K82").Value = 2104.4
$ws.Range("M82").Value = -1743.4
$ws.Range("H85").Value = 2121
$ws.Range("I85").Value = 2104.4
$ws.Range("K85").Value = 2104.4
$ws.Range("M85").Value = -856.4000000000001
$ws.Range("H93").Value = 3168.7273
$ws.Range("I93").Value = 1484.875
$ws.Range("K93").Value = 1484.875
$ws.Range("M93").Value = -236.875
$ws.Range("H113").Value = 9126.5
$ws.Range("I113").Value = 7639.15
$ws.Range("J113").Value = 24000
$ws.Range("K113").Value = 7639.15
$ws.Range("L113").Value = 24000
$ws.Range("M113").Value = -5469.15
$ws.Range("N113").Value = -28340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9825.6
$ws.Range("I81").Value = 3119.5
$ws.Range("K81").Value = 6239
$ws.Range("M81").Value = -5178
$ws.Range("H84").Value = 9825.6
$ws.Range("I84").Value = 3119.5
$ws.Range("K84").Value = 31195
$ws.Range("M84").Value = -25891
$ws.Range("H113").Value = 1524.0834
$ws.Range("I113").Value = 1256.5294
$ws.Range("J113").Value = 2173.8572
$ws.Range("K113").Value = 3769.5882
$ws.Range("L113").Value = 6521.571599999999
$ws.Range("M113").Value = -1599.5882
$ws.Range("N113").Value = -10861.5716
$ws.Range("H115").Value = 95000
$ws.Range("I115").Value = 95000
$ws.Range("K115").Value = 95000
$ws.Range("M115").Value = -93433
$ws.Range("H132").Value = 2861.9048
$ws.Range("I132").Value = 2401.1667
$ws.Range("J132").Value = 3476.2222
$ws.Range("K132").Value = 7203.500100000001
$ws.Range("L132").Value = 10428.6666
$ws.Range("M132").Value = -4673.500100000001
$ws.Range("N132").Value = -15488.6666
